$wb = $excel.ActiveWorkbook

$sheetData = @{
    "ALC" = @(
        @{ Cell = "H19"; Value = 645 },
        @{ Cell = "I19"; Value = 290 },
        @{ Cell = "J19"; Value = 1000 },
        @{ Cell = "K19"; Value = 290 },
        @{ Cell = "L19"; Value = 1000 },
        @{ Cell = "M19"; Value = -115 },
        @{ Cell = "N19"; Value = -1350 },
        @{ Cell = "H116"; Value = 7158.0586 },
        @{ Cell = "J116"; Value = 7419 },
        @{ Cell = "L116"; Value = 7419 },
        @{ Cell = "N116"; Value = -14303 },
        @{ Cell = "H129"; Value = 1507.9474 },
        @{ Cell = "I129"; Value = 972.3333 },
        @{ Cell = "J129"; Value = 2426.1428 },
        @{ Cell = "K129"; Value = 2916.9999 },
        @{ Cell = "L129"; Value = 7278.428400000001 },
        @{ Cell = "M129"; Value = 2083.0001 },
        @{ Cell = "N129"; Value = -17278.4284 },
        @{ Cell = "H132"; Value = 94990 },
        @{ Cell = "I132"; Value = 104309 },
        @{ Cell = "J132"; Value = 1800 },
        @{ Cell = "K132"; Value = 312927 },
        @{ Cell = "L132"; Value = 5400 },
        @{ Cell = "M132"; Value = -310397 },
        @{ Cell = "N132"; Value = -10460 },
        @{ Cell = "H137"; Value = 4547867.5 },
        @{ Cell = "I137"; Value = 1667.2667 },
        @{ Cell = "K137"; Value = 5001.800099999999 },
        @{ Cell = "M137"; Value = -2451.800099999999 },
        @{ Cell = "H138"; Value = 5921.0156 },
        @{ Cell = "I138"; Value = 9230.714 },
        @{ Cell = "J138"; Value = 4341.386 },
        @{ Cell = "K138"; Value = 27692.142 },
        @{ Cell = "L138"; Value = 13024.158 },
        @{ Cell = "M138"; Value = -22552.142 },
        @{ Cell = "N138"; Value = -23304.158 },
        @{ Cell = "H141"; Value = 1890.5834 },
        @{ Cell = "I141"; Value = 1367.6207 },
        @{ Cell = "K141"; Value = 4102.8621 },
        @{ Cell = "M141"; Value = 1077.1379 }
    )
    "ARM" = @(
        @{ Cell = "H28"; Value = 27499.75 },
        @{ Cell = "I28"; Value = 26666.334 },
        @{ Cell = "K28"; Value = 26666.334 },
        @{ Cell = "M28"; Value = -26474.334 },
        @{ Cell = "H92"; Value = 59990.668 },
        @{ Cell = "J92"; Value = 59990.668 },
        @{ Cell = "L92"; Value = 59990.668 },
        @{ Cell = "N92"; Value = -64982.668 },
        @{ Cell = "H93"; Value = 81987.22 },
        @{ Cell = "J93"; Value = 81987.22 },
        @{ Cell = "L93"; Value = 81987.22 },
        @{ Cell = "N93"; Value = -86979.22 },
        @{ Cell = "H95"; Value = 18332.334 },
        @{ Cell = "J95"; Value = 18332.334 },
        @{ Cell = "L95"; Value = 18332.334 },
        @{ Cell = "N95"; Value = -23824.334 },
        @{ Cell = "H96"; Value = 46031.332 },
        @{ Cell = "J96"; Value = 46031.332 },
        @{ Cell = "L96"; Value = 46031.332 },
        @{ Cell = "N96"; Value = -51523.332 },
        @{ Cell = "H97"; Value = 16272.714 },
        @{ Cell = "I97"; Value = 22201.8 },
        @{ Cell = "J97"; Value = 1450 },
        @{ Cell = "K97"; Value = 22201.8 },
        @{ Cell = "L97"; Value = 1450 },
        @{ Cell = "M97"; Value = -21705.8 },
        @{ Cell = "N97"; Value = -2442 },
        @{ Cell = "H99"; Value = 27499.75 },
        @{ Cell = "I99"; Value = 26666.334 },
        @{ Cell = "K99"; Value = 26666.334 },
        @{ Cell = "M99"; Value = -23671.334 },
        @{ Cell = "H101"; Value = 60274.555 },
        @{ Cell = "J101"; Value = 60274.555 },
        @{ Cell = "L101"; Value = 60274.555 },
        @{ Cell = "N101"; Value = -66764.55499999999 },
        @{ Cell = "H103"; Value = 0 },
        @{ Cell = "J103"; Value = 0 },
        @{ Cell = "L103"; Value = 0 },
        @{ Cell = "H104"; Value = 57799.4 },
        @{ Cell = "J104"; Value = 57799.4 },
        @{ Cell = "L104"; Value = 57799.4 },
        @{ Cell = "N104"; Value = -64787.4 },
        @{ Cell = "H105"; Value = 90856.71000000001 },
        @{ Cell = "I105"; Value = 95944.44500000001 },
        @{ Cell = "J105"; Value = 81698.8 },
        @{ Cell = "K105"; Value = 95944.44500000001 },
        @{ Cell = "L105"; Value = 81698.8 },
        @{ Cell = "M105"; Value = -92450.44500000001 },
        @{ Cell = "N105"; Value = -88686.8 },
        @{ Cell = "H106"; Value = 0 },
        @{ Cell = "J106"; Value = 0 },
        @{ Cell = "L106"; Value = 0 },
        @{ Cell = "H132"; Value = 2209.1667 },
        @{ Cell = "I132"; Value = 1263.9166 },
        @{ Cell = "J132"; Value = 4099.6665 },
        @{ Cell = "K132"; Value = 3791.7498 },
        @{ Cell = "L132"; Value = 12298.9995 },
        @{ Cell = "M132"; Value = -1261.7498 },
        @{ Cell = "N132"; Value = -17358.9995 }
    )
    "BSM" = @(
        @{ Cell = "H20"; Value = 1131.7894 },
        @{ Cell = "J20"; Value = 1152.6923 },
        @{ Cell = "L20"; Value = 1152.6923 },
        @{ Cell = "N20"; Value = -1646.6923 },
        @{ Cell = "H105"; Value = 17695.25 },
        @{ Cell = "I105"; Value = 26895.75 },
        @{ Cell = "J105"; Value = 8494.75 },
        @{ Cell = "K105"; Value = 26895.75 },
        @{ Cell = "L105"; Value = 8494.75 },
        @{ Cell = "M105"; Value = -25148.75 },
        @{ Cell = "N105"; Value = -11988.75 }
    )
    "CRP" = @(
        @{ Cell = "H31"; Value = 5216.185 },
        @{ Cell = "I31"; Value = 3335.25 },
        @{ Cell = "J31"; Value = 6008.1577 },
        @{ Cell = "K31"; Value = 3335.25 },
        @{ Cell = "L31"; Value = 6008.1577 },
        @{ Cell = "M31"; Value = -3040.25 },
        @{ Cell = "N31"; Value = -6598.1577 },
        @{ Cell = "H34"; Value = 5216.185 },
        @{ Cell = "I34"; Value = 3335.25 },
        @{ Cell = "J34"; Value = 6008.1577 },
        @{ Cell = "K34"; Value = 3335.25 },
        @{ Cell = "L34"; Value = 6008.1577 },
        @{ Cell = "M34"; Value = -3133.25 },
        @{ Cell = "N34"; Value = -6412.1577 },
        @{ Cell = "H58"; Value = 1945.0244 },
        @{ Cell = "I58"; Value = 1664.48 },
        @{ Cell = "J58"; Value = 2383.375 },
        @{ Cell = "K58"; Value = 1664.48 },
        @{ Cell = "L58"; Value = 2383.375 },
        @{ Cell = "M58"; Value = -1461.48 },
        @{ Cell = "N58"; Value = -2789.375 },
        @{ Cell = "H99"; Value = 1179964.2 },
        @{ Cell = "I99"; Value = 2224926 },
        @{ Cell = "J99"; Value = 4382.5 },
        @{ Cell = "K99"; Value = 2224926 },
        @{ Cell = "L99"; Value = 4382.5 },
        @{ Cell = "M99"; Value = -2223428 },
        @{ Cell = "N99"; Value = -7378.5 },
        @{ Cell = "H126"; Value = 1179964.2 },
        @{ Cell = "I126"; Value = 2224926 },
        @{ Cell = "J126"; Value = 4382.5 },
        @{ Cell = "K126"; Value = 6674778 },
        @{ Cell = "L126"; Value = 13147.5 },
        @{ Cell = "M126"; Value = -6672308 },
        @{ Cell = "N126"; Value = -18087.5 },
        @{ Cell = "H132"; Value = 2357.3962 },
        @{ Cell = "I132"; Value = 2259.2888 },
        @{ Cell = "K132"; Value = 6777.866399999999 },
        @{ Cell = "M132"; Value = -4247.866399999999 },
        @{ Cell = "H136"; Value = 1945.0244 },
        @{ Cell = "I136"; Value = 1664.48 },
        @{ Cell = "J136"; Value = 2383.375 },
        @{ Cell = "K136"; Value = 4993.440000000001 },
        @{ Cell = "L136"; Value = 7150.125 },
        @{ Cell = "M136"; Value = -2443.440000000001 },
        @{ Cell = "N136"; Value = -12250.125 },
        @{ Cell = "H138"; Value = 94365.75 },
        @{ Cell = "J138"; Value = 93987.664 },
        @{ Cell = "L138"; Value = 93987.664 },
        @{ Cell = "N138"; Value = -104267.664 }
    )
    "CUL" = @(
        @{ Cell = "H132"; Value = 1518.75 },
        @{ Cell = "I132"; Value = 1317.7778 },
        @{ Cell = "J132"; Value = 1777.1428 },
        @{ Cell = "K132"; Value = 11860.0002 },
        @{ Cell = "L132"; Value = 15994.2852 },
        @{ Cell = "M132"; Value = -9330.0002 },
        @{ Cell = "N132"; Value = -21054.2852 }
    )
    "GSM" = @(
        @{ Cell = "H15"; Value = 50587 },
        @{ Cell = "J15"; Value = 50587 },
        @{ Cell = "L15"; Value = 50587 },
        @{ Cell = "N15"; Value = -51163 },
        @{ Cell = "H81"; Value = 50587 },
        @{ Cell = "J81"; Value = 50587 },
        @{ Cell = "L81"; Value = 50587 },
        @{ Cell = "N81"; Value = -52583 },
        @{ Cell = "H84"; Value = 50587 },
        @{ Cell = "J84"; Value = 50587 },
        @{ Cell = "L84"; Value = 151761 },
        @{ Cell = "N84"; Value = -161745 },
        @{ Cell = "H97"; Value = 13014.556 },
        @{ Cell = "I97"; Value = 2324.2 },
        @{ Cell = "J97"; Value = 26377.5 },
        @{ Cell = "K97"; Value = 2324.2 },
        @{ Cell = "L97"; Value = 26377.5 },
        @{ Cell = "M97"; Value = -1828.2 },
        @{ Cell = "N97"; Value = -27369.5 },
        @{ Cell = "H133"; Value = 130639.25 },
        @{ Cell = "J133"; Value = 130639.25 },
        @{ Cell = "L133"; Value = 130639.25 },
        @{ Cell = "N133"; Value = -140759.25 }
    )
    "LTW" = @(
        @{ Cell = "H132"; Value = 10488 },
        @{ Cell = "I132"; Value = 3996.182 },
        @{ Cell = "K132"; Value = 11988.546 },
        @{ Cell = "M132"; Value = -9458.545999999998 }
    )
    "WVR" = @(
        @{ Cell = "H81"; Value = 10900.846 },
        @{ Cell = "I81"; Value = 13689.889 },
        @{ Cell = "J81"; Value = 4625.5 },
        @{ Cell = "K81"; Value = 27379.778 },
        @{ Cell = "L81"; Value = 9251 },
        @{ Cell = "M81"; Value = -26318.778 },
        @{ Cell = "N81"; Value = -11373 },
        @{ Cell = "H84"; Value = 10900.846 },
        @{ Cell = "I84"; Value = 13689.889 },
        @{ Cell = "J84"; Value = 4625.5 },
        @{ Cell = "K84"; Value = 136898.89 },
        @{ Cell = "L84"; Value = 46255 },
        @{ Cell = "M84"; Value = -131594.89 },
        @{ Cell = "N84"; Value = -56863 },
        @{ Cell = "H126"; Value = 10200 },
        @{ Cell = "I126"; Value = 16500 },
        @{ Cell = "J126"; Value = 3900 },
        @{ Cell = "K126"; Value = 49500 },
        @{ Cell = "L126"; Value = 11700 },
        @{ Cell = "M126"; Value = -47030 },
        @{ Cell = "N126"; Value = -16640 },
        @{ Cell = "H132"; Value = 30041.055 },
        @{ Cell = "I132"; Value = 38543.63 },
        @{ Cell = "J132"; Value = 4533.3335 },
        @{ Cell = "K132"; Value = 115630.89 },
        @{ Cell = "L132"; Value = 13600.0005 },
        @{ Cell = "M132"; Value = -113100.89 },
        @{ Cell = "N132"; Value = -18660.0005 },
        @{ Cell = "H136"; Value = 26667.195 },
        @{ Cell = "I136"; Value = 35957.586 },
        @{ Cell = "K136"; Value = 107872.758 },
        @{ Cell = "M136"; Value = -105322.758 },
        @{ Cell = "H141"; Value = 150714 },
        @{ Cell = "J141"; Value = 0 },
        @{ Cell = "L141"; Value = 0 }
    )
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($p in $sheetData[$sheetName]) {
        $ws.Range($p.Cell).Value = $p.Value
    }
}

# Cells that are fully removed in the target (no cached value at all)
$clearCells = @(
    @{ Sheet = "ARM"; Cell = "N103" },
    @{ Sheet = "ARM"; Cell = "N106" },
    @{ Sheet = "WVR"; Cell = "N141" }
)

foreach ($c in $clearCells) {
    $ws = $wb.Worksheets.Item($c.Sheet)
    $ws.Range($c.Cell).ClearContents()
}

Write-Output "Applied $($sheetData.Keys.Count) sheets of updates and $($clearCells.Count) clears."